$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix A20: was stored as text ("20072700"), should be a real number
$ws.Range("A20").Value = 20072700

# Add new row 21 with numeric data
$ws.Range("A21").Value = 20072800
$ws.Range("B21").Value = 6398553400
